$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (Price + Volume(1h) columns) to match the latest scrape.
# Price values are stored as plain text (they mix thousands separators with
# decimal points, e.g. "29.190.57"), so we force Text format before writing them
# and then restore the default "Normal" style so no visual/formatting change sticks.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.190.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.859.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7067"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07640"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08415"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.879.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.182"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7085"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.212.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.931"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007810"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.113.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.851"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1586"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.909"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("E30").Value = "  -3.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.393"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.215"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05121"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.8126"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.907"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("E36").Value = "  -2.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.677"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01841"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.698"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.165.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.181"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("E42").Value = "  -1.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.010.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5166"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.773"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.252"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("E51").Value = "  +0.31%  "
